# Auto-generated edit script applying the Exodus_Profits leve-profit value updates
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 932.5
$ws.Range("I18").Value = 932.5
$ws.Range("K18").Value = 932.5
$ws.Range("M18").Value = -648.5
$ws.Range("H33").Value = 174.8
$ws.Range("I33").Value = 187.33333
$ws.Range("J33").Value = 124.666664
$ws.Range("K33").Value = 187.33333
$ws.Range("L33").Value = 124.666664
$ws.Range("M33").Value = 41.66667000000001
$ws.Range("N33").Value = -582.666664
$ws.Range("H92").Value = 1307.5454
$ws.Range("I92").Value = 874.125
$ws.Range("K92").Value = 874.125
$ws.Range("M92").Value = 373.875
$ws.Range("H132").Value = 2292.5454
$ws.Range("I132").Value = 2292.5454
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 6877.6362
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").ClearContents()
$ws.Range("H138").Value = 62502260
$ws.Range("J138").Value = 111113864
$ws.Range("L138").Value = 333341592
$ws.Range("N138").Value = -333351872

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1164.3684
$ws.Range("I2").Value = 965.5294
$ws.Range("J2").Value = 2854.5
$ws.Range("K2").Value = 965.5294
$ws.Range("L2").Value = 2854.5
$ws.Range("M2").Value = -852.5294
$ws.Range("N2").Value = -3080.5
$ws.Range("H32").Value = 6576.255
$ws.Range("I32").Value = 3058.7317
$ws.Range("K32").Value = 3058.7317
$ws.Range("M32").Value = -2771.7317
$ws.Range("H45").Value = 7883137.5
$ws.Range("I45").Value = 2227.5
$ws.Range("J45").Value = 15764048
$ws.Range("K45").Value = 2227.5
$ws.Range("L45").Value = 15764048
$ws.Range("M45").Value = -1850.5
$ws.Range("N45").Value = -15764802
$ws.Range("H46").Value = 29998
$ws.Range("I46").Value = 29998
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 29998
$ws.Range("L46").Value = 0
$ws.Range("M46").ClearContents()
$ws.Range("N46").ClearContents()
$ws.Range("H74").Value = 5753.5293
$ws.Range("I74").Value = 2879.6
$ws.Range("K74").Value = 2879.6
$ws.Range("M74").Value = -2005.6
$ws.Range("H77").Value = 5753.5293
$ws.Range("I77").Value = 2879.6
$ws.Range("K77").Value = 14398
$ws.Range("M77").Value = -10030
$ws.Range("H116").Value = 1164.3684
$ws.Range("I116").Value = 965.5294
$ws.Range("J116").Value = 2854.5
$ws.Range("K116").Value = 965.5294
$ws.Range("L116").Value = 2854.5
$ws.Range("M116").Value = 1328.4706
$ws.Range("N116").Value = -7442.5
$ws.Range("H122").Value = 2086.5334
$ws.Range("I122").Value = 1983.5834
$ws.Range("K122").Value = 5950.7502
$ws.Range("M122").Value = -3500.7502
$ws.Range("H132").Value = 4229.0293
$ws.Range("I132").Value = 3460.6296
$ws.Range("K132").Value = 10381.8888
$ws.Range("M132").Value = -7851.888800000001
$ws.Range("H138").Value = 150000
$ws.Range("J138").Value = 150000
$ws.Range("L138").Value = 150000
$ws.Range("N138").Value = -160280
$ws.Range("H139").Value = 275000
$ws.Range("J139").Value = 275000
$ws.Range("L139").Value = 275000
$ws.Range("N139").Value = -285280

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1164.3684
$ws.Range("I3").Value = 965.5294
$ws.Range("J3").Value = 2854.5
$ws.Range("K3").Value = 965.5294
$ws.Range("L3").Value = 2854.5
$ws.Range("M3").Value = -851.5294
$ws.Range("N3").Value = -3082.5
$ws.Range("H80").Value = 1187.8148
$ws.Range("I80").Value = 837.4
$ws.Range("J80").Value = 1267.4546
$ws.Range("K80").Value = 837.4
$ws.Range("L80").Value = 1267.4546
$ws.Range("M80").Value = 160.6
$ws.Range("N80").Value = -3263.4546
$ws.Range("H83").Value = 1187.8148
$ws.Range("I83").Value = 837.4
$ws.Range("J83").Value = 1267.4546
$ws.Range("K83").Value = 4187
$ws.Range("L83").Value = 6337.273
$ws.Range("M83").Value = 805
$ws.Range("N83").Value = -16321.273
$ws.Range("H86").Value = 6759.7
$ws.Range("I86").Value = 3120
$ws.Range("K86").Value = 3120
$ws.Range("M86").Value = -1997
$ws.Range("H89").Value = 6759.7
$ws.Range("I89").Value = 3120
$ws.Range("K89").Value = 15600
$ws.Range("M89").Value = -9984
$ws.Range("H94").Value = 1221.9333
$ws.Range("I94").Value = 1041
$ws.Range("K94").Value = 1041
$ws.Range("M94").Value = -590

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2106.1292
$ws.Range("I31").Value = 1520.7693
$ws.Range("J31").Value = 2528.889
$ws.Range("K31").Value = 1520.7693
$ws.Range("L31").Value = 2528.889
$ws.Range("M31").Value = -1225.7693
$ws.Range("N31").Value = -3118.889
$ws.Range("H34").Value = 2106.1292
$ws.Range("I34").Value = 1520.7693
$ws.Range("J34").Value = 2528.889
$ws.Range("K34").Value = 1520.7693
$ws.Range("L34").Value = 2528.889
$ws.Range("M34").Value = -1318.7693
$ws.Range("N34").Value = -2932.889
$ws.Range("H97").Value = 23875
$ws.Range("J97").Value = 22285.715
$ws.Range("L97").Value = 22285.715
$ws.Range("N97").Value = -24267.715
$ws.Range("H122").Value = 5113.222
$ws.Range("I122").Value = 4383.8
$ws.Range("K122").Value = 13151.4
$ws.Range("M122").Value = -10701.4
$ws.Range("H132").Value = 1775817.8
$ws.Range("I132").Value = 1519849.1
$ws.Range("K132").Value = 4559547.300000001
$ws.Range("M132").Value = -4557017.300000001

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1277
$ws.Range("I5").Value = 1058.6
$ws.Range("K5").Value = 3175.8
$ws.Range("M5").Value = -3063.8
$ws.Range("H32").Value = 1585.2858
$ws.Range("I32").Value = 192
$ws.Range("K32").Value = 576
$ws.Range("M32").Value = -293
$ws.Range("H34").Value = 1160
$ws.Range("J34").Value = 1400
$ws.Range("L34").Value = 4200
$ws.Range("N34").Value = -4368
$ws.Range("H39").Value = 449
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()
$ws.Range("H48").Value = 3700
$ws.Range("H55").Value = 4921.5
$ws.Range("I55").Value = 337.5
$ws.Range("J55").Value = 6449.5
$ws.Range("K55").Value = 1012.5
$ws.Range("L55").Value = 19348.5
$ws.Range("M55").Value = -835.5
$ws.Range("N55").Value = -19702.5
$ws.Range("H121").Value = 1002868.8
$ws.Range("I121").Value = 3241.3333
$ws.Range("J121").Value = 1431280.6
$ws.Range("K121").Value = 9723.999899999999
$ws.Range("L121").Value = 4293841.800000001
$ws.Range("M121").Value = -8413.999899999999
$ws.Range("N121").Value = -4296461.800000001
$ws.Range("H131").Value = 1674.5
$ws.Range("J131").Value = 2138.4285
$ws.Range("L131").Value = 6415.2855
$ws.Range("N131").Value = -16495.2855
$ws.Range("H135").Value = 1277
$ws.Range("I135").Value = 1058.6
$ws.Range("K135").Value = 9527.4
$ws.Range("M135").Value = -6992.4

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 5330.4473
$ws.Range("I132").Value = 2819.0908
$ws.Range("K132").Value = 8457.2724
$ws.Range("M132").Value = -5927.2724

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 14107.8
$ws.Range("J7").Value = 4599.8335
$ws.Range("L7").Value = 4599.8335
$ws.Range("N7").Value = -4823.8335
$ws.Range("H22").Value = 75973.11
$ws.Range("I22").Value = 2488.2
$ws.Range("K22").Value = 2488.2
$ws.Range("M22").Value = -2193.2
$ws.Range("H25").Value = 11499.929
$ws.Range("J25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("N25").ClearContents()
$ws.Range("H27").Value = 75973.11
$ws.Range("I27").Value = 2488.2
$ws.Range("K27").Value = 2488.2
$ws.Range("M27").Value = -2381.2
$ws.Range("H126").Value = 14107.8
$ws.Range("J126").Value = 4599.8335
$ws.Range("L126").Value = 13799.5005
$ws.Range("N126").Value = -18739.5005
$ws.Range("H136").Value = 5390.4287
$ws.Range("I136").Value = 8011.8335
$ws.Range("K136").Value = 24035.5005
$ws.Range("M136").Value = -21485.5005

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H6").Value = 2750
$ws.Range("J6").Value = 2500
$ws.Range("L6").Value = 2500
$ws.Range("N6").Value = -2730
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("M11").ClearContents()
$ws.Range("H122").Value = 2601.2122
$ws.Range("I122").Value = 2771.24
$ws.Range("J122").Value = 2069.875
$ws.Range("K122").Value = 8313.719999999999
$ws.Range("L122").Value = 6209.625
$ws.Range("M122").Value = -5863.719999999999
$ws.Range("N122").Value = -11109.625
$ws.Range("H136").Value = 2116.5386
$ws.Range("I136").Value = 2116.5386
$ws.Range("K136").Value = 6349.6158
$ws.Range("M136").Value = -3799.6158
